$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: D(date serial), J(Volumen), K(Precio minimo), L(Precio maximo),
# M(Precio promedio ponderado), N(Unidad de comercializacion), O(Origen),
# P(Precio $/Kg), Q(Kg o Unidades) -- for sheet rows 2..100
$data = @(
    @(44594,400,31000,32000,31500,"$/malla 25 kilos","Provincia del Elquí",1260,25),
    @(44610,400,30000,31000,30500,"$/malla 25 kilos","Provincia del Elquí",1220,25),
    @(44558,400,30000,31000,30500,"$/malla 25 kilos","Provincia de Limarí",1220,25),
    @(44330,520,32000,33000,32500,"$/malla 25 kilos","Provincia del Elquí",1300,25),
    @(44333,400,31500,32000,31750,"$/malla 25 kilos","Provincia del Elquí",1270,25),
    @(44316,600,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44309,600,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44650,480,27000,28000,27500,"$/malla 25 kilos","Provincia del Elquí",1100,25),
    @(44596,500,31000,32000,31500,"$/malla 25 kilos","Provincia del Elquí",1260,25),
    @(44321,400,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44638,400,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44277,560,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44251,700,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44680,360,27000,28000,27500,"$/malla 25 kilos","Provincia del Elquí",1100,25),
    @(44335,400,30000,31000,30500,"$/malla 25 kilos","Provincia de Limarí",1220,25),
    @(44295,600,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44216,600,36000,37000,36500,"$/malla 25 kilos","Provincia del Elquí",1460,25),
    @(44246,600,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44291,500,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44657,500,27000,28000,27500,"$/malla 25 kilos","Provincia del Elquí",1100,25),
    @(44244,640,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44307,560,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44230,500,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44586,600,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44690,480,24500,25000,24750,"$/malla 25 kilos","Provincia del Elquí",990,25),
    @(44649,400,27000,28000,27500,"$/malla 25 kilos","Provincia del Elquí",1100,25),
    @(44641,500,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44588,500,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44223,660,32500,33000,32750,"$/malla 25 kilos","Provincia del Elquí",1310,25),
    @(44566,400,16000,17000,16500,"$/caja 15 kilos","Provincia de Limarí",1100,15),
    @(44566,600,31000,32000,31500,"$/malla 25 kilos","Provincia del Elquí",1260,25),
    @(44678,400,26000,27000,26500,"$/malla 25 kilos","Provincia del Elquí",1060,25),
    @(44642,400,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44265,760,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44601,300,31000,32000,31500,"$/malla 25 kilos","Provincia del Elquí",1260,25),
    @(44200,400,31000,32000,31500,"$/malla 25 kilos","Provincia de Limarí",1260,25),
    @(44187,400,37000,38000,37500,"$/malla 25 kilos","Provincia de Limarí",1500,25),
    @(44237,600,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44636,500,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44607,600,30000,31000,30500,"$/malla 25 kilos","Provincia del Elquí",1220,25),
    @(44643,540,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44568,700,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44670,600,26000,27000,26500,"$/malla 25 kilos","Provincia del Elquí",1060,25),
    @(44218,400,34000,35000,34500,"$/malla 25 kilos","Provincia del Elquí",1380,25),
    @(44571,600,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44225,600,31000,32000,31500,"$/malla 25 kilos","Provincia de Limarí",1260,25),
    @(44663,560,24000,25000,24500,"$/malla 25 kilos","Provincia del Elquí",980,25),
    @(44587,400,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44239,600,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44600,400,30000,31000,30500,"$/malla 25 kilos","Provincia del Elquí",1220,25),
    @(44274,600,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44573,400,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44671,540,26000,27000,26500,"$/malla 25 kilos","Provincia del Elquí",1060,25),
    @(44272,600,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44293,500,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44231,500,25000,26000,25500,"$/malla 25 kilos","Provincia del Elquí",1020,25),
    @(44252,520,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44202,600,29000,30000,29500,"$/malla 25 kilos","Provincia de Limarí",1180,25),
    @(44648,480,27000,28000,27500,"$/malla 25 kilos","Provincia del Elquí",1100,25),
    @(44629,400,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44253,660,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44281,640,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44279,560,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44615,520,31000,32000,31500,"$/malla 25 kilos","Provincia del Elquí",1260,25),
    @(44300,400,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44692,400,25000,26000,25500,"$/malla 25 kilos","Provincia del Elquí",1020,25),
    @(44298,500,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44603,520,31000,32000,31500,"$/malla 25 kilos","Provincia del Elquí",1260,25),
    @(44685,440,26000,27000,26500,"$/malla 25 kilos","Provincia del Elquí",1060,25),
    @(44687,400,25000,26000,25500,"$/malla 25 kilos","Provincia del Elquí",1020,25),
    @(44302,600,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44662,460,24000,25000,24500,"$/malla 25 kilos","Provincia del Elquí",980,25),
    @(44655,440,27000,28000,27500,"$/malla 25 kilos","Provincia del Elquí",1100,25),
    @(44582,500,30000,31000,30500,"$/malla 25 kilos","Provincia del Elquí",1220,25),
    @(44631,520,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44575,400,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44580,600,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44614,540,31000,32000,31500,"$/malla 25 kilos","Provincia del Elquí",1260,25),
    @(44260,600,27000,28000,27500,"$/malla 25 kilos","Provincia del Elquí",1100,25),
    @(44589,560,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44557,400,30000,31000,30500,"$/malla 25 kilos","Provincia de Limarí",1220,25),
    @(44323,600,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44221,460,35000,36000,35500,"$/malla 25 kilos","Provincia del Elquí",1420,25),
    @(44278,400,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44245,540,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44552,400,35000,36000,35500,"$/malla 25 kilos","Provincia de Limarí",1420,25),
    @(44559,360,30000,31000,30500,"$/malla 25 kilos","Provincia de Limarí",1220,25),
    @(44559,400,32000,33000,32500,"$/malla 25 kilos","Provincia del Elquí",1300,25),
    @(44664,500,26000,27000,26500,"$/malla 25 kilos","Provincia del Elquí",1060,25),
    @(44656,400,27000,28000,27500,"$/malla 25 kilos","Provincia del Elquí",1100,25),
    @(44628,560,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44608,500,30000,31000,30500,"$/malla 25 kilos","Provincia del Elquí",1220,25),
    @(44209,600,36000,37000,36500,"$/malla 25 kilos","Provincia del Elquí",1460,25),
    @(44659,400,27000,28000,27500,"$/malla 25 kilos","Provincia del Elquí",1100,25),
    @(44505,300,37000,38000,37500,"$/malla 25 kilos","Perú",1500,25),
    @(44326,500,29000,30000,29500,"$/malla 25 kilos","Provincia del Elquí",1180,25),
    @(44238,520,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44264,400,28000,29000,28500,"$/malla 25 kilos","Provincia del Elquí",1140,25),
    @(44547,200,22000,22500,22250,"$/caja 15 kilos","Provincia de Limarí",1483,15)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 4).Value  = $row[0]   # D - Fecha
    $ws.Cells.Item($r, 10).Value = $row[1]   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $row[2]   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $row[3]   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $row[4]   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $row[5]   # N - Unidad de comercializacion
    $ws.Cells.Item($r, 15).Value = $row[6]   # O - Origen
    $ws.Cells.Item($r, 16).Value = $row[7]   # P - Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $row[8]   # Q - Kg o Unidades
}
